$p = $ppt.ActivePresentation
Write-Output "Before count: $($p.Designs.Count)"
try {
    $d2 = $p.Designs.Add("{12345678-1234-1234-1234-123456789012}", 2)
    Write-Output "Added: $d2"
} catch {
    Write-Output "ERR add: $_"
}
Write-Output "After count: $($p.Designs.Count)"
